# Repull data, push all data, mean calculation
# Update column F (dSF) values on the active sheet's data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 4
    3  = 2
    4  = -2
    5  = -5
    6  = -1
    7  = -1
    8  = -5
    9  = -4
    10 = -2
    11 = -1
    12 = 1
    13 = 2
    15 = 0
    17 = -1
    19 = -3
    20 = 0
    21 = -1
    22 = 3
    24 = 6
    25 = -3
    26 = -1
    27 = -3
    29 = -2
    30 = 5
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
